# Sexting intensity escalation fix for AshleyJourney sheet (PPV0 phase detection bug)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AshleyJourney")

$ws.Range("B4").Value  = "watch me let go... this is only for you 😊"
$ws.Range("B6").Value  = "cum with me... I'm letting go right now, watch me"
$ws.Range("B7").Value  = "I'm so close... I can feel it building everywhere and I don't want to cum alone 😊"
$ws.Range("B8").Value  = "oh god I can't hold on 💕"
$ws.Range("B9").Value  = "this is what you made me do and you need to see every second of it 😊"
$ws.Range("B11").Value = "I need you to watch what I'm doing right now... you have to see this"
$ws.Range("B12").Value = "my fingers keep going deeper and faster and my whole body is shaking 😊"
$ws.Range("B13").Value = "I'm touching my pussy and imagining it's your hands on me... I need more"
$ws.Range("B14").Value = "fuck 💕"
$ws.Range("B15").Value = "look what you did to me... I couldn't stop 😊"
$ws.Range("B17").Value = "tell me what you want me to do to myself right now... I'll do anything you say"
$ws.Range("B18").Value = "I'm so wet right now babe... you have no idea what your words do to my body 😊"
$ws.Range("B19").Value = "but I can't stop now... my fingers are between my legs and it's all because of you"
$ws.Range("B20").Value = "oh god... I can't believe I just did that 💕"
$ws.Range("B21").Value = "I want you to see what you're doing to me right now 😊"
$ws.Range("B23").Value = "my hand keeps sliding lower and I can't stop it babe... I don't even want to"
$ws.Range("B24").Value = "I keep running my hands down my body and everything is so sensitive... it's like every touch is amplified because of you"
$ws.Range("B25").Value = "you really liked that? knowing you saw me is making my heart race so fast right now 💕"
